$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled data: update column F (dSF) values for the affected rows
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = 0
$ws.Range("F17").Value = -1
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -2
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = -4
$ws.Range("F26").Value = -7
$ws.Range("F30").Value = -1
